# Remove the two leading empty paragraphs that immediately follow the
# first table in the document (they carried the placeholder spacing
# used before the new "R Development Best Practices" content existed).
#
# Both paragraphs are completely empty (just a paragraph mark), so we
# delete them as single-character ranges taken right at the end of the
# table. Deleting one character at a time (rather than a single two
# character range) reliably merges away the *leading* paragraph and
# keeps the one after it, matching how Word actually collapses empty
# paragraphs when their marks are deleted.

$d = $word.ActiveDocument

$table = $d.Tables.Item(1)
$afterTable = $table.Range.End

$d.Range($afterTable, $afterTable + 1).Delete()
$d.Range($afterTable, $afterTable + 1).Delete()
